$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of language/value pairs (rows 2-21), sorted descending by value,
# with Bengali and Uzbek removed.
$data = @(
    @("English", 24.3804824838931),
    @("Chinese", 12.50612318325456),
    @("Spanish", 7.440664061081822),
    @("Japanese", 5.602459728343781),
    @("Arabic", 5.058356011939543),
    @("German", 4.822664839451362),
    @("Russian", 3.738306842148356),
    @("Portuguese", 3.49491804240869),
    @("French", 3.015734838146286),
    @("Malay-Indonesian", 2.693108611207126),
    @("Italian", 2.654654056454913),
    @("Korean", 1.659161428273157),
    @("Persian", 1.587627847122411),
    @("Turkish", 1.475335024953981),
    @("Dutch", 1.42485829133396),
    @("Thai", 0.9782416845756626),
    @("Polish", 0.8737857907173976),
    @("Urdu", 0.8636404448634035),
    @("Vietnamese", 0.4764158394806041),
    @("Swedish", 0.4751894926604449)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the now-unused rows 22 and 23 (previously held Uzbek and Vietnamese).
$ws.Rows.Item(22).Delete()
$ws.Rows.Item(22).Delete()
